# CF1: Divide adjacent repeat rules. Fix instruments
# The "No bass change" / "Bass change" labels are renamed to
# "No cantus change" / "Cantus change" (instrument terminology fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "No cantus change"
$ws.Range("G3").Value = "Cantus change"

# Reflect the author's final active selection in the saved file.
$ws.Range("G4").Select() | Out-Null
